$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the whole "Meta description: ..." paragraph (2nd
#    paragraph of the document) - its content is being relocated to
#    the end of the document (see step 2 below).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)

if ($metaPara.Range.Text -notmatch "Meta description") {
    # Defensive fallback - locate the paragraph by its text instead of
    # a hard-coded index, in case the document shape differs slightly.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match "Meta description") {
            $metaPara = $p
            break
        }
    }
}

# Keep a copy of its formatted content (rich text incl. bold run) on
# the clipboard before deleting the paragraph - it will be pasted back
# in right before the closing "Prompt:" paragraph.
$metaPara.Range.Copy()
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new paragraph right before the last paragraph
#    ("Prompt: Create a feature image ...") containing the bold text
#    "Play Extra Wild Slot Free Today".
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.Paste()

$newParaIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$newRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Play Extra Wild Slot Free Today", 2) | Out-Null

$newRange2 = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newRange2.Find.Execute(": Read our Extra Wild slot review and play it for free. Enjoy unique Wild symbol feature with hidden multiplier and original symbols. Get started now!", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph - the old image
#    generation prompt is swapped out for the meta-description copy.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Find.Execute("Prompt: Create a feature image for Extra Wild that follows the game's luxury and gemstone theme. The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing on a treasure chest overflowing with precious gems, gold coins, and diamond-encrusted watches. The warrior should be holding a large diamond in one hand and flashing a winning grin, indicating the success players can achieve in the game. The background should have a blue velvet texture, similar to that of the game grid, to tie in the game's visual aspect. Overall, the image should evoke feelings of luxury, wealth, and success, fitting for a game about precious gemstones.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our Extra Wild slot review and play it for free. Enjoy unique Wild symbol feature with hidden multiplier and original symbols. Get started now!", 2) | Out-Null

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
